# dev_train.xlsx update
# - Rename sheet "Set2_7" -> "Set2_7_THISONE" (this is the sheet the author is
#   currently working with / uploading).
# - Change the active tab / selection state:
#     * "Set7_6" (previously the active tab, with selection A1:J6) loses focus
#       and its remembered selection becomes K8.
#     * "Set2_7_THISONE" becomes the new active tab, with remembered
#       selection M2:M5 (active cell M2).

$wb = $excel.ActiveWorkbook

# --- Update the selection remembered on the sheet that is losing focus ---
$wsOld = $wb.Worksheets.Item("Set7_6")
$wsOld.Range("K8").Select() | Out-Null

# --- Switch to / activate the sheet that becomes the new active tab ---
$wsNew = $wb.Worksheets.Item("Set2_7")
$wsNew.Activate() | Out-Null
$wsNew.Range("M2:M5").Select() | Out-Null

# --- Rename the newly active sheet ---
$wsNew.Name = "Set2_7_THISONE"
